$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '70.023.73'
$ws.Range('E2').Value = '  -0.49%  '
$ws.Range('D3').Value = '3.541.30'
$ws.Range('E3').Value = '  -0.16%  '
$ws.Range('E4').Value = '  -0.08%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '603.87'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -2.06%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '197.23'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +5.75%  '
$ws.Range('E7').Value = '  -0.43%  '
$ws.Range('E8').Value = '  -0.04%  '
$ws.Range('E9').Value = '  -3.27%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.653'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -1.03%  '
$ws.Range('E12').Value = '  -2.09%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '9.53'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -1.33%  '
$ws.Range('D14').Value = '4.104.15'
$ws.Range('E14').Value = '  -0.15%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '602.12'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -3.10%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '19.27'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +0.78%  '
$ws.Range('B17').Value = 'Uniswap'
$ws.Range('C17').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '12.83'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -0.24%  '
$ws.Range('B18').Value = 'WrappedBTC'
$ws.Range('C18').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D18').Value = '70.150.62'
$ws.Range('E18').Value = '  -0.31%  '
$ws.Range('D19').Value = '3.541.40'
$ws.Range('E19').Value = '  -0.50%  '
$ws.Range('E20').Value = '  +0.38%  '
$ws.Range('E21').Value = '  -0.41%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '17.88'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +1.67%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '5.26'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +3.66%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '102.31'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -2.26%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '4.62'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -2.38%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '3.13'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +3.03%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '10.95'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -0.39%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '9.59'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -2.57%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '33.66'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -1.81%  '
$ws.Range('E30').Value = '  +0.59%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '4.33'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +17.94%  '
$ws.Range('E32').Value = '  +0.74%  '
$ws.Range('E33').Value = '  -1.44%  '
$ws.Range('D35').Value = '0.0₃0839'
$ws.Range('E35').Value = '  +7.33%  '
$ws.Range('D36').Value = '3.781.42'
$ws.Range('E36').Value = '  +6.60%  '
$ws.Range('E37').Value = '  +0.12%  '
$ws.Range('E38').Value = '  -3.62%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '3.65'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +1.03%  '
$ws.Range('E40').Value = '  -1.93%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '36.67'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -1.70%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '495.30'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -7.43%  '
$ws.Range('E43').Value = '  -3.12%  '
$ws.Range('E44').Value = '  -2.93%  '
$ws.Range('E46').Value = '  -1.70%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '3.30'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -2.41%  '
$ws.Range('E48').Value = '  +0.23%  '
$ws.Range('E49').Value = '  -4.34%  '
$ws.Range('E50').Value = '  +2.59%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '129.66'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -4.03%  '
